$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.60235125631333
$ws.Range("C2").Value = 10.35297868524203
$ws.Range("D2").Value = 14.86559037952263
$ws.Range("E2").Value = 16.26969319116602
$ws.Range("G2").Value = 29.58962610198869
$ws.Range("H2").Value = 14.48003855791578
$ws.Range("J2").Value = 9.275767282234055
$ws.Range("O2").Value = 22.15713680413623
$ws.Range("B3").Value = 15.88034851117316
$ws.Range("C3").Value = 9.709045822229147
$ws.Range("D3").Value = 14.80044589252489
$ws.Range("E3").Value = 16.204637478889
$ws.Range("G3").Value = 29.71485210227291
$ws.Range("H3").Value = 14.55259742355067
$ws.Range("J3").Value = 9.283182266756789
$ws.Range("O3").Value = 22.27594940005278
$ws.Range("B4").Value = 15.42018450494558
$ws.Range("C4").Value = 9.289982135954709
$ws.Range("D4").Value = 14.76377113487245
$ws.Range("E4").Value = 16.16836745281891
$ws.Range("G4").Value = 29.80619064815159
$ws.Range("H4").Value = 14.60053837602079
$ws.Range("J4").Value = 9.289174943423262
$ws.Range("O4").Value = 22.35597753025775
$ws.Range("B5").Value = 15.22866073336323
$ws.Range("C5").Value = 9.113287375673362
$ws.Range("D5").Value = 14.74967304201776
$ws.Range("E5").Value = 16.15452181654948
$ws.Range("G5").Value = 29.84701254892002
$ws.Range("H5").Value = 14.62092536776583
$ws.Range("J5").Value = 9.291979185297263
$ws.Range("O5").Value = 22.39036023304108
$ws.Range("B6").Value = 15.1966240624933
$ws.Range("C6").Value = 9.083590813103621
$ws.Range("D6").Value = 14.7473835380891
$ws.Range("E6").Value = 16.15227950218622
$ws.Range("G6").Value = 29.85400747895702
$ws.Range("H6").Value = 14.62436194236819
$ws.Range("J6").Value = 9.29246670389985
$ws.Range("O6").Value = 22.39617613249474
$ws.Range("B7").Value = 15.41761742373273
$ws.Range("C7").Value = 9.287623098154969
$ws.Range("D7").Value = 14.76357755887602
$ws.Range("E7").Value = 16.16817692790972
$ws.Range("G7").Value = 29.80672665203599
$ws.Range("H7").Value = 14.60080988001125
$ws.Range("J7").Value = 9.289211295877802
$ws.Range("O7").Value = 22.35643407022457
$ws.Range("B8").Value = 16.35703285637407
$ws.Range("C8").Value = 10.13588329342839
$ws.Range("D8").Value = 14.84244571897594
$ws.Range("E8").Value = 16.24650614350473
$ws.Range("G8").Value = 29.62978596344811
$ws.Range("H8").Value = 14.50435235971665
$ws.Range("J8").Value = 9.278025221446795
$ws.Range("O8").Value = 22.19662927739881
$ws.Range("B9").Value = 18.05679156224383
$ws.Range("C9").Value = 11.61043484536726
$ws.Range("D9").Value = 15.02293572258256
$ws.Range("E9").Value = 16.42873953441852
$ws.Range("G9").Value = 29.39882474257864
$ws.Range("H9").Value = 14.34216840777776
$ws.Range("J9").Value = 9.267506459396198
$ws.Range("O9").Value = 21.93979857646943
$ws.Range("B10").Value = 19.2090373370403
$ws.Range("C10").Value = 12.57793861869174
$ws.Range("D10").Value = 15.17048571283122
$ws.Range("E10").Value = 16.57930221079355
$ws.Range("G10").Value = 29.30159477628506
$ws.Range("H10").Value = 14.23955603662917
$ws.Range("J10").Value = 9.26672499660785
$ws.Range("O10").Value = 21.78611337907333
$ws.Range("B11").Value = 19.71068137225912
$ws.Range("C11").Value = 12.99289636838321
$ws.Range("D11").Value = 15.24066352466681
$ws.Range("E11").Value = 16.65123027841781
$ws.Range("G11").Value = 29.27341702199147
$ws.Range("H11").Value = 14.1964902927683
$ws.Range("J11").Value = 9.267873214806551
$ws.Range("O11").Value = 21.72391299547737
$ws.Range("B12").Value = 19.89729826838039
$ws.Range("C12").Value = 13.14641340842233
$ws.Range("D12").Value = 15.26766007652591
$ws.Range("E12").Value = 16.67894437618771
$ws.Range("G12").Value = 29.26507761954688
$ws.Range("H12").Value = 14.18070399807773
$ws.Range("J12").Value = 9.268523659542907
$ws.Range("O12").Value = 21.70147748908715
$ws.Range("B13").Value = 19.85725725893407
$ws.Range("C13").Value = 13.11351166237903
$ws.Range("D13").Value = 15.2618274206076
$ws.Range("E13").Value = 16.67295474691603
$ws.Range("G13").Value = 29.26676966270128
$ws.Range("H13").Value = 14.18408061788902
$ws.Range("J13").Value = 9.26837399450565
$ws.Range("O13").Value = 21.70625948486583
$ws.Range("B14").Value = 19.7261020416233
$ws.Range("C14").Value = 13.00559881456781
$ws.Range("D14").Value = 15.24287619963212
$ws.Range("E14").Value = 16.65350088825805
$ws.Range("G14").Value = 29.27268410049446
$ws.Range("H14").Value = 14.1951810750634
$ws.Range("J14").Value = 9.267922409299059
$ws.Range("O14").Value = 21.72204474406246
$ws.Range("B15").Value = 19.6453272192621
$ws.Range("C15").Value = 12.93902790976254
$ws.Range("D15").Value = 15.23132241100751
$ws.Range("E15").Value = 16.64164635412604
$ws.Range("G15").Value = 29.27661103765976
$ws.Range("H15").Value = 14.20204844112305
$ws.Range("J15").Value = 9.267673863654373
$ws.Range("O15").Value = 21.73185960264911
$ws.Range("B16").Value = 19.17579107001951
$ws.Range("C16").Value = 12.55031371882471
$ws.Range("D16").Value = 15.16595943070067
$ws.Range("E16").Value = 16.57466920228886
$ws.Range("G16").Value = 29.30376127076956
$ws.Range("H16").Value = 14.24244337740519
$ws.Range("J16").Value = 9.266680167994366
$ws.Range("O16").Value = 21.79033431505535
$ws.Range("B17").Value = 18.88189678053897
$ws.Range("C17").Value = 12.30540396539689
$ws.Range("D17").Value = 15.12663181370431
$ws.Range("E17").Value = 16.53444874787722
$ws.Range("G17").Value = 29.32454545621051
$ws.Range("H17").Value = 14.26815133070465
$ws.Range("J17").Value = 9.266455314412472
$ws.Range("O17").Value = 21.82818862053081
$ws.Range("B18").Value = 18.71074219947267
$ws.Range("C18").Value = 12.16216949737228
$ws.Range("D18").Value = 15.10430042613813
$ws.Range("E18").Value = 16.51163956614088
$ws.Range("G18").Value = 29.33800974418333
$ws.Range("H18").Value = 14.28327781721629
$ws.Range("J18").Value = 9.26646756259259
$ws.Range("O18").Value = 21.85068672605332
$ws.Range("B19").Value = 18.65243259734409
$ws.Range("C19").Value = 12.11326596331446
$ws.Range("D19").Value = 15.09678953600529
$ws.Range("E19").Value = 16.50397303578733
$ws.Range("G19").Value = 29.34282709314828
$ws.Range("H19").Value = 14.28845770613803
$ws.Range("J19").Value = 9.266496042182894
$ws.Range("O19").Value = 21.8584284796169
$ws.Range("B20").Value = 18.91340198586823
$ws.Range("C20").Value = 12.33172020326385
$ws.Range("D20").Value = 15.13078854057267
$ws.Range("E20").Value = 16.53869681823806
$ws.Range("G20").Value = 29.32217651587633
$ws.Range("H20").Value = 14.26537947557972
$ws.Range("J20").Value = 9.266464601637733
$ws.Range("O20").Value = 21.82408382901892
$ws.Range("B21").Value = 19.76471705044098
$ws.Range("C21").Value = 13.03739362044277
$ws.Range("D21").Value = 15.24843133238955
$ws.Range("E21").Value = 16.65920217090296
$ws.Range("G21").Value = 29.27088346106335
$ws.Range("H21").Value = 14.19190642335976
$ws.Range("J21").Value = 9.268049203645258
$ws.Range("O21").Value = 21.71737780206912
$ws.Range("B22").Value = 20.30156673213235
$ws.Range("C22").Value = 13.4775030649493
$ws.Range("D22").Value = 15.32776713237515
$ws.Range("E22").Value = 16.74072802065554
$ws.Range("G22").Value = 29.25095519305937
$ws.Range("H22").Value = 14.1469298442617
$ws.Range("J22").Value = 9.270341398673118
$ws.Range("O22").Value = 21.6541626012149
$ws.Range("B23").Value = 20.01685820813953
$ws.Range("C23").Value = 13.24453715848919
$ws.Range("D23").Value = 15.28520607140657
$ws.Range("E23").Value = 16.69696885547621
$ws.Range("G23").Value = 29.26034065674179
$ws.Range("H23").Value = 14.17065558327365
$ws.Range("J23").Value = 9.269003256732647
$ws.Range("O23").Value = 21.68730176021074
$ws.Range("B24").Value = 18.89916529825119
$ws.Range("C24").Value = 12.31983020916678
$ws.Range("D24").Value = 15.1289084152624
$ws.Range("E24").Value = 16.53677528594582
$ws.Range("G24").Value = 29.32324279576127
$ws.Range("H24").Value = 14.26663155163164
$ws.Range("J24").Value = 9.266459962018008
$ws.Range("O24").Value = 21.82593731592339
$ws.Range("B25").Value = 17.61331370505355
$ws.Range("C25").Value = 11.23184793400959
$ws.Range("D25").Value = 14.97142194316469
$ws.Range("E25").Value = 16.37645148708101
$ws.Range("G25").Value = 29.44869798086624
$ws.Range("H25").Value = 14.38314572948062
$ws.Range("J25").Value = 9.269130806334404
$ws.Range("O25").Value = 22.00316874956335
